{"js": "// Update the date line and the 25 division-problem cells in the table.\n// Each old value is unique in the document, so a plain text search +\n// whole-match replace is unambiguous for every entry.\nconst replacements = [\n  [\"2025-05-04 Sunday\", \"2025-05-05 Monday\"],\n  [\"977\u00f79=\", \"403\u00f73=\"],\n  [\"596\u00f75=\", \"599\u00f78=\"],\n  [\"238\u00f72=\", \"275\u00f72=\"],\n  [\"982\u00f78=\", \"269\u00f79=\"],\n  [\"644\u00f76=\", \"541\u00f76=\"],\n  [\"166\u00f76=\", \"795\u00f78=\"],\n  [\"172\u00f73=\", \"545\u00f75=\"],\n  [\"561\u00f79=\", \"405\u00f72=\"],\n  [\"501\u00f77=\", \"783\u00f78=\"],\n  [\"911\u00f75=\", \"393\u00f74=\"],\n  [\"224\u00f74=\", \"880\u00f79=\"],\n  [\"130\u00f75=\", \"346\u00f74=\"],\n  [\"263\u00f77=\", \"879\u00f78=\"],\n  [\"370\u00f78=\", \"439\u00f73=\"],\n  [\"602\u00f77=\", \"937\u00f73=\"],\n  [\"539\u00f75=\", \"341\u00f73=\"],\n  [\"972\u00f79=\", \"897\u00f77=\"],\n  [\"994\u00f74=\", \"580\u00f75=\"],\n  [\"465\u00f75=\", \"451\u00f74=\"],\n  [\"297\u00f72=\", \"834\u00f76=\"],\n  [\"766\u00f78=\", \"516\u00f77=\"],\n  [\"177\u00f74=\", \"219\u00f73=\"],\n  [\"411\u00f77=\", \"936\u00f76=\"],\n  [\"649\u00f77=\", \"968\u00f76=\"],\n  [\"494\u00f72=\", \"504\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-04 Sunday\", \"2025-05-05 Monday\"),\n    @(\"977\u00f79=\", \"403\u00f73=\"),\n    @(\"596\u00f75=\", \"599\u00f78=\"),\n    @(\"238\u00f72=\", \"275\u00f72=\"),\n    @(\"982\u00f78=\", \"269\u00f79=\"),\n    @(\"644\u00f76=\", \"541\u00f76=\"),\n    @(\"166\u00f76=\", \"795\u00f78=\"),\n    @(\"172\u00f73=\", \"545\u00f75=\"),\n    @(\"561\u00f79=\", \"405\u00f72=\"),\n    @(\"501\u00f77=\", \"783\u00f78=\"),\n    @(\"911\u00f75=\", \"393\u00f74=\"),\n    @(\"224\u00f74=\", \"880\u00f79=\"),\n    @(\"130\u00f75=\", \"346\u00f74=\"),\n    @(\"263\u00f77=\", \"879\u00f78=\"),\n    @(\"370\u00f78=\", \"439\u00f73=\"),\n    @(\"602\u00f77=\", \"937\u00f73=\"),\n    @(\"539\u00f75=\", \"341\u00f73=\"),\n    @(\"972\u00f79=\", \"897\u00f77=\"),\n    @(\"994\u00f74=\", \"580\u00f75=\"),\n    @(\"465\u00f75=\", \"451\u00f74=\"),\n    @(\"297\u00f72=\", \"834\u00f76=\"),\n    @(\"766\u00f78=\", \"516\u00f77=\"),\n    @(\"177\u00f74=\", \"219\u00f73=\"),\n    @(\"411\u00f77=\", \"936\u00f76=\"),\n    @(\"649\u00f77=\", \"968\u00f76=\"),\n    @(\"494\u00f72=\", \"504\u00f78=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $found = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
